$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1224.1562
$ws.Range("I70").Value = 1276.6666
$ws.Range("K70").Value = 3829.9998
$ws.Range("M70").Value = -3559.9998

$ws.Range("H73").Value = 1224.1562
$ws.Range("I73").Value = 1276.6666
$ws.Range("K73").Value = 3829.9998
$ws.Range("M73").Value = -2893.9998

$ws.Range("H113").Value = 3198.4736
$ws.Range("I113").Value = 2719.5
$ws.Range("J113").Value = 3730.6667
$ws.Range("K113").Value = 2719.5
$ws.Range("L113").Value = 3730.6667
$ws.Range("M113").Value = 534.5
$ws.Range("N113").Value = -10238.6667

$ws.Range("H116").Value = 3838.5715
$ws.Range("I116").Value = 3812.818
$ws.Range("J116").Value = 3933
$ws.Range("K116").Value = 3812.818
$ws.Range("L116").Value = 3933
$ws.Range("M116").Value = -370.8180000000002
$ws.Range("N116").Value = -10817

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9130.621999999999
$ws.Range("I32").Value = 8227.985000000001
$ws.Range("J32").Value = 17770.143
$ws.Range("K32").Value = 8227.985000000001
$ws.Range("L32").Value = 17770.143
$ws.Range("M32").Value = -7940.985000000001
$ws.Range("N32").Value = -18344.143

$ws.Range("H45").Value = 2011.5294
$ws.Range("I45").Value = 1882
$ws.Range("J45").Value = 2102.2
$ws.Range("K45").Value = 1882
$ws.Range("L45").Value = 2102.2
$ws.Range("M45").Value = -1505
$ws.Range("N45").Value = -2856.2

$ws.Range("H61").Value = 3628.0833
$ws.Range("I61").Value = 4440.375
$ws.Range("K61").Value = 4440.375
$ws.Range("M61").Value = -4228.375

$ws.Range("H63").Value = 5287.846
$ws.Range("I63").Value = 3416.6667
$ws.Range("J63").Value = 6891.7144
$ws.Range("K63").Value = 3416.6667
$ws.Range("L63").Value = 6891.7144
$ws.Range("M63").Value = -2730.6667
$ws.Range("N63").Value = -8263.714400000001

$ws.Range("H66").Value = 5287.846
$ws.Range("I66").Value = 3416.6667
$ws.Range("J66").Value = 6891.7144
$ws.Range("K66").Value = 17083.3335
$ws.Range("L66").Value = 34458.572
$ws.Range("M66").Value = -13651.3335
$ws.Range("N66").Value = -41322.572

$ws.Range("H101").Value = 30602
$ws.Range("J101").Value = 30602
$ws.Range("L101").Value = 30602
$ws.Range("N101").Value = -37092

$ws.Range("H110").Value = 1597
$ws.Range("I110").Value = 1397.5333
$ws.Range("J110").Value = 2594.3333
$ws.Range("K110").Value = 1397.5333
$ws.Range("L110").Value = 2594.3333
$ws.Range("M110").Value = 647.4666999999999
$ws.Range("N110").Value = -6684.3333

$ws.Range("H119").Value = 24339.4
$ws.Range("J119").Value = 24339.4
$ws.Range("L119").Value = 24339.4
$ws.Range("N119").Value = -34015.4

$ws.Range("H132").Value = 854820.5600000001
$ws.Range("I132").Value = 1703609.5
$ws.Range("J132").Value = 6031.6787
$ws.Range("K132").Value = 5110828.5
$ws.Range("L132").Value = 18095.0361
$ws.Range("M132").Value = -5108298.5
$ws.Range("N132").Value = -23155.0361

$ws.Range("H136").Value = 3628.0833
$ws.Range("I136").Value = 4440.375
$ws.Range("K136").Value = 13321.125
$ws.Range("M136").Value = -10771.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5537.846
$ws.Range("I134").Value = 2327.682
$ws.Range("J134").Value = 9692.177
$ws.Range("K134").Value = 6983.045999999999
$ws.Range("L134").Value = 29076.531
$ws.Range("M134").Value = -4448.045999999999
$ws.Range("N134").Value = -34146.531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1295.6316
$ws.Range("I107").Value = 252.54546
$ws.Range("J107").Value = 2729.875
$ws.Range("K107").Value = 252.54546
$ws.Range("L107").Value = 2729.875
$ws.Range("M107").Value = 1667.45454
$ws.Range("N107").Value = -6569.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 630.1818
$ws.Range("I6").Value = 198
$ws.Range("K6").Value = 594
$ws.Range("M6").Value = -481

$ws.Range("H17").Value = 820.2222
$ws.Range("J17").Value = 1664
$ws.Range("L17").Value = 4992
$ws.Range("N17").Value = -5330

$ws.Range("H22").Value = 2726.9285
$ws.Range("J22").Value = 2775.5557
$ws.Range("L22").Value = 8326.667099999999
$ws.Range("N22").Value = -8664.667099999999

$ws.Range("H25").Value = 47623104
$ws.Range("I25").Value = 380
$ws.Range("J25").Value = 55560224
$ws.Range("K25").Value = 1140
$ws.Range("L25").Value = 166680672
$ws.Range("M25").Value = -971
$ws.Range("N25").Value = -166681010

$ws.Range("H27").Value = 2726.9285
$ws.Range("J27").Value = 2775.5557
$ws.Range("L27").Value = 8326.667099999999
$ws.Range("N27").Value = -8530.667099999999

$ws.Range("H30").Value = 47623104
$ws.Range("I30").Value = 380
$ws.Range("J30").Value = 55560224
$ws.Range("K30").Value = 1140
$ws.Range("L30").Value = 166680672
$ws.Range("M30").Value = -1038
$ws.Range("N30").Value = -166680876

$ws.Range("H34").Value = 1475.5
$ws.Range("I34").Value = 140
$ws.Range("J34").Value = 2429.4285
$ws.Range("K34").Value = 420
$ws.Range("L34").Value = 7288.2855
$ws.Range("M34").Value = -336
$ws.Range("N34").Value = -7456.2855

$ws.Range("H39").Value = 957.0909
$ws.Range("J39").Value = 909.75
$ws.Range("L39").Value = 2729.25
$ws.Range("N39").Value = -3317.25

$ws.Range("H55").Value = 1746.25
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 2161.6667
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 6485.000100000001
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -6839.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 10494.625
$ws.Range("I99").Value = 7112.1816
$ws.Range("K99").Value = 7112.1816
$ws.Range("M99").Value = -4866.1816

$ws.Range("H109").Value = 20275
$ws.Range("J109").Value = 20275
$ws.Range("L109").Value = 20275
$ws.Range("N109").Value = -22355

$ws.Range("H132").Value = 2086406.5
$ws.Range("I132").Value = 5954812.5
$ws.Range("J132").Value = 3418.4614
$ws.Range("K132").Value = 17864437.5
$ws.Range("L132").Value = 10255.3842
$ws.Range("M132").Value = -17861907.5
$ws.Range("N132").Value = -15315.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 38050.6
$ws.Range("I132").Value = 71587.664
$ws.Range("J132").Value = 4513.533
$ws.Range("K132").Value = 214762.992
$ws.Range("L132").Value = 13540.599
$ws.Range("M132").Value = -212232.992
$ws.Range("N132").Value = -18600.599

$ws.Range("H136").Value = 1755.5518
$ws.Range("I136").Value = 1450.9166
$ws.Range("J136").Value = 1970.5883
$ws.Range("K136").Value = 4352.7498
$ws.Range("L136").Value = 5911.7649
$ws.Range("M136").Value = -1802.7498
$ws.Range("N136").Value = -11011.7649

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2478.6875
$ws.Range("I126").Value = 2843.476
$ws.Range("J126").Value = 1782.2727
$ws.Range("K126").Value = 8530.428
$ws.Range("L126").Value = 5346.8181
$ws.Range("M126").Value = -6060.428
$ws.Range("N126").Value = -10286.8181

$ws.Range("H132").Value = 1795.0677
$ws.Range("I132").Value = 1586.091
$ws.Range("J132").Value = 2060.3076
$ws.Range("K132").Value = 4758.272999999999
$ws.Range("L132").Value = 6180.9228
$ws.Range("M132").Value = -2228.272999999999
$ws.Range("N132").Value = -11240.9228
